# Expand the "User Database" sheet with a new data row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 / column A gains the value "edr" (dimension grows from A1:C1 to A1:C2).
$ws.Range("A2").Value = "edr"
